$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add row 48 ---------------------------------------------------------
# Clone the cell formatting from row 44 (same alternating colour-band style
# family: date cell + "content" cells, with the last/"empty" variant style
# used for the Remarque/Reference columns) so the new rows keep the
# journal's existing look.
$fmtCols = @(1, 2, 3, 4, 5, 7)   # A, B, C, D, E, G
foreach ($col in $fmtCols) {
    $ws.Cells.Item(44, $col).Copy()
    $ws.Cells.Item(48, $col).PasteSpecial(-4122)  # xlPasteFormats
}
$ws.Cells.Item(44, 8).Copy()                      # H (empty-cell variant)
$ws.Cells.Item(48, 6).PasteSpecial(-4122)         # F48 uses that same style
$ws.Cells.Item(48, 8).PasteSpecial(-4122)         # H48
$excel.CutCopyMode = $false

$ws.Range("A48").Value = 44267
$ws.Range("B48").Value = "Documentation"
$ws.Range("C48").Value = "6h"
$ws.Range("D48").Value = "Mise a jour du rapport de travail"
$ws.Range("E48").Value = "Oui"
$ws.Range("G48").Value = "Non"
$ws.Rows.Item(48).RowHeight = 30

# --- Add row 49 ---------------------------------------------------------
foreach ($col in $fmtCols) {
    $ws.Cells.Item(44, $col).Copy()
    $ws.Cells.Item(49, $col).PasteSpecial(-4122)
}
$ws.Cells.Item(44, 8).Copy()
$ws.Cells.Item(49, 6).PasteSpecial(-4122)
$ws.Cells.Item(49, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A49").Value = 44267
$ws.Range("B49").Value = "React native (Frontend)"
$ws.Range("C49").Value = "2h"
$ws.Range("D49").Value = "Test de l'application"
$ws.Range("E49").Value = "Oui"
$ws.Range("G49").Value = "Non"
$ws.Rows.Item(49).RowHeight = 31.5

# --- View state: zoom out a bit and land the selection on the new row --
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("F49").Select()
